$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tableau recherche")
$ws.Range("A148").Value = "hello"
Write-Host "step1 ok"
$v = $ws.Range("A148").Value
Write-Host "step2 ok"
Write-Host $v
Write-Host "step3 ok"
